$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("leaderboard2")
$ws.Range("D8").Value = 1064
$ws.Range("G9").Value = 709
$ws.Range("D10").Value = 964
$ws.Range("J10").Value = 378
$ws.Range("C12").Value = "Brybry_"
$ws.Range("D12").Value = 948
$ws.Range("J12").Value = 325
$ws.Range("C13").Value = "Nikof_"
$ws.Range("B14").Value = "Dernière update le 13.03.25 à 17:29"

$ws = $wb.Worksheets.Item("leaderboard3")
$ws.Range("F3").Value = "Mickaplow"
$ws.Range("G3").Value = 189
$ws.Range("I3").Value = "Kaatsup"
$ws.Range("L3").Value = "JLTootmy"
$ws.Range("M3").Value = 23
$ws.Range("F4").Value = "Brybry_"
$ws.Range("I4").Value = "AntoineDaniel_"
$ws.Range("L4").Value = "Grimkujow"
$ws.Range("M4").Value = 18
$ws.Range("F5").Value = "_Linca"
$ws.Range("G5").Value = 166
$ws.Range("L5").Value = "KennyStream"
$ws.Range("M5").Value = 17
$ws.Range("D6").Value = 337
$ws.Range("I6").Value = "CrocodyleTV"
$ws.Range("J6").Value = 41
$ws.Range("L6").Value = "Angle_Droit"
$ws.Range("M6").Value = 16
$ws.Range("I7").Value = "LadySundae"
$ws.Range("J7").Value = 40
$ws.Range("L7").Value = "XoTrixy"
$ws.Range("M7").Value = 13
$ws.Range("D8").Value = 332
$ws.Range("L8").Value = "Hiro_Ammar"
$ws.Range("M8").Value = 10
$ws.Range("I9").Value = "Onutrem"
$ws.Range("J9").Value = 30
$ws.Range("L9").Value = "Maxouzboub"
$ws.Range("M9").Value = 7
$ws.Range("G10").Value = 122
$ws.Range("I10").Value = "KyriaaTV"
$ws.Range("J10").Value = 29
$ws.Range("I11").Value = "nemenems"
$ws.Range("J11").Value = 29
$ws.Range("I12").Value = "LuttiLutti"
$ws.Range("J12").Value = 28
$ws.Range("C13").Value = "DFG_DrFeelgood"
$ws.Range("I13").Value = "Terraciid"
$ws.Range("J13").Value = 28
$ws.Range("B14").Value = "Dernière update le 13.03.25 à 17:29"

$ws = $wb.Worksheets.Item("leaderboard4")
$ws.Range("J8").Value = 18
$ws.Range("D9").Value = 72
$ws.Range("B14").Value = "Dernière update le 13.03.25 à 17:29"

$ws = $wb.Worksheets.Item("leaderboard5")
$ws.Range("F3").Value = "TheGuill84"
$ws.Range("G3").Value = 20
$ws.Range("I3").Value = "HexakiI"
$ws.Range("O3").Value = "CrocodyleTV"
$ws.Range("F4").Value = "MoMaN_uS"
$ws.Range("G4").Value = 19
$ws.Range("I4").Value = "ChloeRamdani"
$ws.Range("O4").Value = "Maxouzboub"
$ws.Range("D5").Value = 358
$ws.Range("F5").Value = "_Linca"
$ws.Range("I5").Value = "Onutrem"
$ws.Range("J5").Value = 6
$ws.Range("O5").Value = "HarryLafranc"
$ws.Range("D6").Value = 269
$ws.Range("F6").Value = "LadySundae"
$ws.Range("G6").Value = 13
$ws.Range("L6").Value = "XoTrixy"
$ws.Range("O6").Value = "Grimkujow"
$ws.Range("D7").Value = 182
$ws.Range("F7").Value = "Theorus_"
$ws.Range("G7").Value = 12
$ws.Range("I7").Value = "JimmyBoyyy"
$ws.Range("J7").Value = 5
$ws.Range("L7").Value = "Pepito_kawazakii"
$ws.Range("O7").Value = "Horty_"
$ws.Range("F8").Value = "LuttiLutti"
$ws.Range("G8").Value = 9
$ws.Range("L8").Value = "Mynth0s"
$ws.Range("C9").Value = "Mickaplow"
$ws.Range("D9").Value = 85
$ws.Range("I9").Value = "Etoiles"
$ws.Range("J9").Value = 4
$ws.Range("L9").Value = "ZeratoR"
$ws.Range("C10").Value = "SakorRos"
$ws.Range("D10").Value = 76
$ws.Range("F10").Value = "nemenems"
$ws.Range("G10").Value = 8
$ws.Range("I10").Value = "KyriaaTV"
$ws.Range("L10").Value = "Wingobear"
$ws.Range("D11").Value = 54
$ws.Range("F11").Value = "Gom4rt_"
$ws.Range("I11").Value = "Bytell2"
$ws.Range("C12").Value = "AntoineDaniel_"
$ws.Range("D12").Value = 30
$ws.Range("I12").Value = "Angle_Droit"
$ws.Range("L12").Value = "Hiro_Ammar"
$ws.Range("C13").Value = "Nikof_"
$ws.Range("D13").Value = 24
$ws.Range("F13").Value = "gobgg"
$ws.Range("G13").Value = 7
$ws.Range("I13").Value = "Elspawn"
$ws.Range("J13").Value = 4
$ws.Range("L13").Value = "BagheraJones"
$ws.Range("B14").Value = "Dernière update le 13.03.25 à 17:29"
